$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Table Filter" row (row 22)
$ws.Range("A22").Value = "Работа по созданию функционала загрузки данных (Групповые кнопки полей/фильтров)"
$ws.Range("B22").Value = 2
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C22").Value = 43566

# Update selection to match the after-state
$ws.Range("A16").Select()
